$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)
Write-Host "slide1 HasNotesPage: [$($s1.HasNotesPage)]"
$new = $s1.Duplicate()
$new.MoveTo(3)
Write-Host "new HasNotesPage: [$($new.HasNotesPage)]"
$np = $new.NotesPage
Write-Host "new HasNotesPage after access: [$($new.HasNotesPage)]"
Write-Host "np shapes: $($np.Shapes.Count)"
